$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.2095808383233533
$ws.Range("C2").Value = 0.5149700598802395
$ws.Range("J2").Value = 0.01796407185628742
$ws.Range("P2").Value = 0.1377245508982036
$ws.Range("S2").Value = 0.1197604790419162
$ws.Range("B3").Value = 0.005714285714285714
$ws.Range("C3").Value = 0.02857142857142857
$ws.Range("J3").Value = 0.01714285714285714
$ws.Range("P3").Value = 0.7142857142857143
$ws.Range("S3").Value = 0.2342857142857143
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2708333333333333
$ws.Range("B6").Value = 0.04621848739495799
$ws.Range("D6").Value = 0.02100840336134454
$ws.Range("F6").Value = 0.02941176470588235
$ws.Range("J6").Value = 0.2184873949579832
$ws.Range("O6").Value = 0.02100840336134454
$ws.Range("Q6").Value = 0.1974789915966386
$ws.Range("R6").Value = 0.06302521008403361
$ws.Range("S6").Value = 0.4033613445378151
$ws.Range("B7").Value = 0.09183673469387756
$ws.Range("D7").Value = 0.03061224489795918
$ws.Range("F7").Value = 0.06122448979591837
$ws.Range("J7").Value = 0.1122448979591837
$ws.Range("O7").Value = 0.02551020408163265
$ws.Range("Q7").Value = 0.1530612244897959
$ws.Range("R7").Value = 0.08673469387755102
$ws.Range("S7").Value = 0.4387755102040816
$ws.Range("B8").Value = 0.08021390374331551
$ws.Range("D8").Value = 0.0071301247771836
$ws.Range("F8").Value = 0.07843137254901961
$ws.Range("J8").Value = 0.09269162210338681
$ws.Range("O8").Value = 0.0124777183600713
$ws.Range("Q8").Value = 0.1925133689839572
$ws.Range("R8").Value = 0.06060606060606061
$ws.Range("S8").Value = 0.4759358288770054
$ws.Range("B9").Value = 0.1172839506172839
$ws.Range("D9").Value = 0.0308641975308642
$ws.Range("F9").Value = 0.05555555555555555
$ws.Range("J9").Value = 0.1049382716049383
$ws.Range("O9").Value = 0.006172839506172839
$ws.Range("Q9").Value = 0.1728395061728395
$ws.Range("R9").Value = 0.07407407407407407
$ws.Range("S9").Value = 0.4382716049382716
$ws.Range("B10").Value = 0.1177304964539007
$ws.Range("D10").Value = 0.02127659574468085
$ws.Range("E10").Value = 0.001418439716312057
$ws.Range("F10").Value = 0.07872340425531915
$ws.Range("J10").Value = 0.1092198581560284
$ws.Range("O10").Value = 0.01418439716312057
$ws.Range("Q10").Value = 0.2205673758865248
$ws.Range("R10").Value = 0.06808510638297872
$ws.Range("S10").Value = 0.3687943262411347
$ws.Range("G11").Value = 0.1698717948717949
$ws.Range("J11").Value = 0.08333333333333333
$ws.Range("K11").Value = 0.2147435897435897
$ws.Range("L11").Value = 0.5160256410256411
$ws.Range("S11").Value = 0.01602564102564102
$ws.Range("G12").Value = 0.7926829268292683
$ws.Range("J12").Value = 0.1768292682926829
$ws.Range("K12").Value = 0.01219512195121951
$ws.Range("L12").Value = 0.01219512195121951
$ws.Range("S12").Value = 0.006097560975609756
$ws.Range("G13").Value = 0.5952380952380952
$ws.Range("J13").Value = 0.3571428571428572
$ws.Range("S13").Value = 0.04761904761904762
$ws.Range("F15").Value = 0.01716738197424893
$ws.Range("H15").Value = 0.2317596566523605
$ws.Range("I15").Value = 0.03433476394849785
$ws.Range("J15").Value = 0.3733905579399142
$ws.Range("K15").Value = 0.06008583690987124
$ws.Range("M15").Value = 0.008583690987124463
$ws.Range("O15").Value = 0.06437768240343347
$ws.Range("S15").Value = 0.2103004291845494
$ws.Range("F16").Value = 0.0303030303030303
$ws.Range("H16").Value = 0.2222222222222222
$ws.Range("I16").Value = 0.0505050505050505
$ws.Range("J16").Value = 0.4393939393939394
$ws.Range("K16").Value = 0.08080808080808081
$ws.Range("M16").Value = 0.005050505050505051
$ws.Range("O16").Value = 0.05555555555555555
$ws.Range("S16").Value = 0.1161616161616162
$ws.Range("F17").Value = 0.007692307692307693
$ws.Range("H17").Value = 0.1865384615384615
$ws.Range("I17").Value = 0.08269230769230769
$ws.Range("J17").Value = 0.4711538461538461
$ws.Range("K17").Value = 0.08846153846153847
$ws.Range("M17").Value = 0.01153846153846154
$ws.Range("O17").Value = 0.05576923076923077
$ws.Range("S17").Value = 0.09615384615384616
$ws.Range("F18").Value = 0.01162790697674419
$ws.Range("H18").Value = 0.2093023255813954
$ws.Range("I18").Value = 0.08139534883720931
$ws.Range("J18").Value = 0.4593023255813953
$ws.Range("K18").Value = 0.0872093023255814
$ws.Range("M18").Value = 0.01162790697674419
$ws.Range("O18").Value = 0.05813953488372093
$ws.Range("S18").Value = 0.08139534883720931
$ws.Range("F19").Value = 0.01079136690647482
$ws.Range("H19").Value = 0.2366906474820144
$ws.Range("I19").Value = 0.06330935251798561
$ws.Range("J19").Value = 0.3956834532374101
$ws.Range("K19").Value = 0.1064748201438849
$ws.Range("M19").Value = 0.02302158273381295
$ws.Range("N19").Value = 0.0007194244604316547
$ws.Range("O19").Value = 0.06978417266187051
$ws.Range("S19").Value = 0.09352517985611511
